# Commit message: "added try except to excel write function"
# The underlying Python script re-ran its statsmodels summary export and
# (per the xlsx diff) only the embedded "Date:"/"Time:" stamps inside each
# of the three OLS-summary cells changed - the rest of each summary text
# is identical. Replace those two timestamp substrings, on every sheet's
# B2 summary cell, in place (same-length replacement keeps the
# fixed-width column alignment of the statsmodels text table intact).

$wb = $excel.ActiveWorkbook

$oldDate = "Sun, 29 Dec 2019"
$newDate = "Wed, 01 Jan 2020"
$oldTime = "16:11:11"
$newTime = "23:18:49"

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("B2")
    $text = [string]$cell.Value()
    if ($text -and ($text.Contains($oldDate) -or $text.Contains($oldTime))) {
        $text = $text.Replace($oldDate, $newDate)
        $text = $text.Replace($oldTime, $newTime)
        $cell.Value = $text
    }
}
